# "update markerset to logtrack properties"
#
# The "properties" sheet lists, per model class, the UI-editable property
# rows. This edit:
#   1. Renames the `logtrack` class to `d3-logtrack` (rows 25-37, formerly
#      25-39).
#   2. Inserts a new `markerSet` property row for `d3-logtrack` (new row 38),
#      pushing the existing `zoneSet` row down to 39 and `showZoneSet` down
#      to 40. The shared "Zonation" section header becomes
#      "Zonation & MarkerSet" for both the zoneSet and markerSet rows.
#   3. Adds a reference-spec column (G/H, headed by new column "refSpec")
#      used by curve.endDepth / curve.startDepth to point at
#      well.bottomDepth / well.topDepth respectively.
#   4. Adds new `zoneset` and `zone` classes with their own property rows
#      (rows 41-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")
$ws.Activate()

# --- 1. New "refSpec" header column (G) ------------------------------------
$ws.Range("G1").Value = "refSpec"

# --- 2. curve.endDepth / curve.startDepth gain a wiref refSpec -----------
$ws.Range("F9").Value  = "wiref"
$ws.Range("G9").Value  = "well.bottomDepth"
$ws.Range("H9").Value  = "well"

$ws.Range("F10").Value = "wiref"
$ws.Range("G10").Value = "well.topDepth"
$ws.Range("H10").Value = "well"

# --- 3. logtrack -> d3-logtrack -------------------------------------------
$logtrackRows = 25..38
foreach ($r in $logtrackRows) {
    $ws.Range("A$r").Value = "d3-logtrack"
}

# --- 4. Insert the new markerSet row at 38, shifting zoneSet/showZoneSet --
$ws.Range("B38").Value = "markerSet"
$ws.Range("C38").Value = "MarkerSet"
$ws.Range("D38").Value = "use "
$ws.Range("E38").Value = "Zonation & MarkerSet"
$ws.Range("F38").Value = "wimarkerset"

$ws.Range("A39").Value = "d3-logtrack"
$ws.Range("B39").Value = "zoneSet"
$ws.Range("C39").Value = "ZoneSet"
$ws.Range("D39").Value = "use "
$ws.Range("E39").Value = "Zonation & MarkerSet"
$ws.Range("F39").Value = "wizoneset"

$ws.Range("A40").Value = "d3-logtrack"
$ws.Range("B40").Value = "showZoneSet"
$ws.Range("C40").ClearContents()
$ws.Range("D40").Value = "notuse "
$ws.Range("E40").ClearContents()
$ws.Range("F40").ClearContents()

# --- 5. New "zoneset" class -------------------------------------------------
$ws.Range("A41").Value = "zoneset"
$ws.Range("B41").Value = "name"
$ws.Range("C41").Value = "Name"
$ws.Range("D41").Value = "use "
$ws.Range("E41").Value = "Properties"

# --- 6. New "zone" class -----------------------------------------------------
$ws.Range("A42").Value = "zone"
$ws.Range("B42").Value = "name"
$ws.Range("C42").Value = "Name"
$ws.Range("D42").Value = "use "
$ws.Range("E42").Value = "Header"

$ws.Range("A43").Value = "zone"
$ws.Range("B43").Value = "endDepth"
$ws.Range("C43").Value = "End Depth"
$ws.Range("D43").Value = "readonly"
$ws.Range("E43").Value = "Depths"

$ws.Range("A44").Value = "zone"
$ws.Range("B44").Value = "startDepth"
$ws.Range("C44").Value = "Start Depth"
$ws.Range("D44").Value = "readonly"
$ws.Range("E44").Value = "Depths"

$ws.Range("A45").Value = "zone"
$ws.Range("B45").Value = "showName"
$ws.Range("C45").Value = "Show Name"
$ws.Range("D45").Value = "use "
$ws.Range("E45").Value = "Options"
$ws.Range("F45").Value = "checkbox"

$ws.Range("A46").Value = "zone"
$ws.Range("B46").Value = "showOnTrack"
$ws.Range("C46").Value = "Show On Track"
$ws.Range("D46").Value = "use "
$ws.Range("E46").Value = "Options"
$ws.Range("F46").Value = "checkbox"

# --- 7. Cosmetics: column widths + selection (best effort) ----------------
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(7).AutoFit()

$ws.Range("H40").Select()

Write-Output "applied logtrack/markerset edit"
